# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Reorder two pairs of country rows (sharedStrings text swap) ---
# Santa Lucia <-> Timor Oriental (rows 204 / 205)
$ws.Range("A204").Value = "Timor Oriental"
$ws.Range("A205").Value = "Santa Lucia"

# Montserrat <-> Islas Malvinas (rows 214 / 215), including their updated stats
$ws.Range("A214").Value = "Islas Malvinas"
$ws.Range("B214").Value = 13
$ws.Range("C214").Value = 0
$ws.Range("D214").Value = 13
$ws.Range("E214").Value = 0
$ws.Range("F214").Value = 0
$ws.Range("G214").Value = 0
$ws.Range("H214").Value = 0

$ws.Range("A215").Value = "Montserrat"
$ws.Range("B215").Value = 13
$ws.Range("C215").Value = 0
$ws.Range("D215").Value = 12
$ws.Range("E215").Value = 0
$ws.Range("F215").Value = 0
$ws.Range("G215").Value = 0
$ws.Range("H215").Value = 1

# --- Updated "last refreshed" timestamp string ---
$ws.Range("A1").Value = "Datos actualizados a 21 de Septiembre de 2020 a las 07:31"

# --- Updated country statistics ---
# Row 5 - India
$ws.Range("B5").Value = 5487580
$ws.Range("C5").Value = 1968
$ws.Range("D5").Value = 4396399
$ws.Range("E5").Value = 1003272

# Row 27 - Israel
$ws.Range("B27").Value = 188427
$ws.Range("C27").Value = 525
$ws.Range("D27").Value = 135991
$ws.Range("E27").Value = 51180

# Row 59 - Uzbekistan
$ws.Range("B59").Value = 51789
$ws.Range("C59").Value = 149
$ws.Range("E59").Value = 3422
$ws.Range("G59").Value = 2
$ws.Range("H59").Value = 435

# Row 66 - Kirguistan
$ws.Range("B66").Value = 45416
$ws.Range("D66").Value = 41578
$ws.Range("E66").Value = 2775

# Row 77 - El Salvador
$ws.Range("D77").Value = 21575
$ws.Range("E77").Value = 5166
$ws.Range("G77").Value = 1
$ws.Range("H77").Value = 812

# Row 136 - Tailandia
$ws.Range("D136").Value = 3342
$ws.Range("E136").Value = 105
